$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dados")

# Remove the "Id" column (column A) entirely - Nome/Telefone shift left.
$ws.Range("A1").EntireColumn.Delete()

# Update phone numbers (content changed as part of the edit).
$ws.Range("B3").Value = "+553196449238"
$ws.Range("B2").Value = "+553183349238"

# Nome (now column A) gets a wider, explicit custom width.
$ws.Range("A1").EntireColumn.ColumnWidth = 20.1

# Selection / active cell per the diff.
$ws.Range("B2").Select()
